$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.877.95'
$ws.Range('E2').Value = '  -2.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.753.77'
$ws.Range('E3').Value = '  -4.61%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.26'
$ws.Range('E5').Value = '  -8.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5094'
$ws.Range('E7').Value = '  -5.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.37'
$ws.Range('E8').Value = '  -5.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2791'
$ws.Range('E9').Value = '  -5.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06215'
$ws.Range('E10').Value = '  -10.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.744.95'
$ws.Range('E11').Value = '  -5.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.79'
$ws.Range('E12').Value = '  -9.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.06968'
$ws.Range('E13').Value = '  -3.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6142'
$ws.Range('E14').Value = '  -15.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.528'
$ws.Range('E15').Value = '  -9.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.63'
$ws.Range('E16').Value = '  -12.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9997'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.883.28'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006964'
$ws.Range('E20').Value = '  -11.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.68'
$ws.Range('E21').Value = '  -15.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.969.09'
$ws.Range('E22').Value = '  -5.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.090'
$ws.Range('E23').Value = '  -10.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.273'
$ws.Range('E24').Value = '  -12.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.238'
$ws.Range('E25').Value = '  -10.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.79'
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.489'
$ws.Range('E27').Value = '  -12.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.825'
$ws.Range('E28').Value = '  -15.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.06'
$ws.Range('E29').Value = '  -11.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '103.81'
$ws.Range('E30').Value = '  -6.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08221'
$ws.Range('E31').Value = '  -7.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.713'
$ws.Range('E32').Value = '  -12.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.507'
$ws.Range('E33').Value = '  -13.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04534'
$ws.Range('E34').Value = '  -6.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9987'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.632'
$ws.Range('E36').Value = '  -9.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9929'
$ws.Range('E37').Value = '  -12.10%  '
$ws.Range('E38').Value = '  -15.07%  '
$ws.Range('E39').Value = '  -12.79%  '
$ws.Range('E40').Value = '  -8.85%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '103.94'
$ws.Range('E41').Value = '  -2.95%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.905'
$ws.Range('E43').Value = '  -16.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3884'
$ws.Range('E44').Value = '  -17.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7451'
$ws.Range('E45').Value = '  -17.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.932'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05412'
$ws.Range('E47').Value = '  -6.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1116'
$ws.Range('E48').Value = '  -10.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.017'
$ws.Range('E49').Value = '  -18.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.26'
$ws.Range('E50').Value = '  -12.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.92'
$ws.Range('E51').Value = '  -11.90%  '
